# Scheduled runner update: refresh Leve-profit market price columns
# (currentAveragePrice / NQ / HQ and derived profit columns) across the
# per-profession sheets with newly scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1023
$ws.Range("I18").Value = 864
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 864
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -580
$ws.Range("N18").Value = -2068
$ws.Range("H114").Value = 31129.143
$ws.Range("J114").Value = 31129.143
$ws.Range("L114").Value = 31129.143
$ws.Range("N114").Value = -39807.143
$ws.Range("H124").Value = 52477
$ws.Range("J124").Value = 52477
$ws.Range("L124").Value = 52477
$ws.Range("N124").Value = -62297

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 43857
$ws.Range("J111").Value = 43857
$ws.Range("L111").Value = 43857
$ws.Range("N111").Value = -52037
$ws.Range("H113").Value = 36160.168
$ws.Range("J113").Value = 36160.168
$ws.Range("L113").Value = 36160.168
$ws.Range("N113").Value = -44838.168
$ws.Range("H114").Value = 31046.5
$ws.Range("J114").Value = 31046.5
$ws.Range("L114").Value = 31046.5
$ws.Range("N114").Value = -39724.5
$ws.Range("H118").Value = 49998
$ws.Range("J118").Value = 49998
$ws.Range("L118").Value = 49998
$ws.Range("N118").Value = -53312
$ws.Range("H119").Value = 51690
$ws.Range("J119").Value = 51690
$ws.Range("L119").Value = 51690
$ws.Range("N119").Value = -61366
$ws.Range("H121").Value = 31739.75
$ws.Range("J121").Value = 31739.75
$ws.Range("L121").Value = 31739.75
$ws.Range("N121").Value = -35233.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 433.66666
$ws.Range("I22").Value = 400.5
$ws.Range("K22").Value = 400.5
$ws.Range("M22").Value = -227.5
$ws.Range("H108").Value = 46663.25
$ws.Range("J108").Value = 46663.25
$ws.Range("L108").Value = 46663.25
$ws.Range("N108").Value = -54343.25
$ws.Range("H111").Value = 35925.5
$ws.Range("J111").Value = 35925.5
$ws.Range("L111").Value = 35925.5
$ws.Range("N111").Value = -44105.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49780
$ws.Range("J20").Value = 49780
$ws.Range("L20").Value = 49780
$ws.Range("N20").Value = -50252
$ws.Range("H30").Value = 49780
$ws.Range("J30").Value = 49780
$ws.Range("L30").Value = 49780
$ws.Range("N30").Value = -49962
$ws.Range("H31").Value = 8700.888999999999
$ws.Range("I31").Value = 2426.3333
$ws.Range("J31").Value = 14975.444
$ws.Range("K31").Value = 2426.3333
$ws.Range("L31").Value = 14975.444
$ws.Range("M31").Value = -2131.3333
$ws.Range("N31").Value = -15565.444
$ws.Range("H34").Value = 8700.888999999999
$ws.Range("I34").Value = 2426.3333
$ws.Range("J34").Value = 14975.444
$ws.Range("K34").Value = 2426.3333
$ws.Range("L34").Value = 14975.444
$ws.Range("M34").Value = -2224.3333
$ws.Range("N34").Value = -15379.444
$ws.Range("H60").Value = 23602.777
$ws.Range("J60").Value = 23602.777
$ws.Range("L60").Value = 23602.777
$ws.Range("N60").Value = -24624.777
$ws.Range("H110").Value = 40985.668
$ws.Range("J110").Value = 40985.668
$ws.Range("L110").Value = 40985.668
$ws.Range("N110").Value = -49165.668
$ws.Range("H112").Value = 40492
$ws.Range("J112").Value = 40492
$ws.Range("L112").Value = 40492
$ws.Range("N112").Value = -43446
$ws.Range("H116").Value = 44970
$ws.Range("J116").Value = 44970
$ws.Range("L116").Value = 44970
$ws.Range("N116").Value = -54148
$ws.Range("H119").Value = 48753
$ws.Range("J119").Value = 48753
$ws.Range("L119").Value = 48753
$ws.Range("N119").Value = -58429
$ws.Range("H128").Value = 49780
$ws.Range("J128").Value = 49780
$ws.Range("L128").Value = 49780
$ws.Range("N128").Value = -59740

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 339078.78
$ws.Range("I125").Value = 3000000
$ws.Range("J125").Value = 6463.625
$ws.Range("K125").Value = 9000000
$ws.Range("L125").Value = 19390.875
$ws.Range("M125").Value = -8995080
$ws.Range("N125").Value = -29230.875
$ws.Range("H126").Value = 252962.75
$ws.Range("I126").Value = 800566
$ws.Range("J126").Value = 4052.182
$ws.Range("K126").Value = 2401698
$ws.Range("L126").Value = 12156.546
$ws.Range("M126").Value = -2396758
$ws.Range("N126").Value = -22036.546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 48694
$ws.Range("J110").Value = 48694
$ws.Range("L110").Value = 48694
$ws.Range("N110").Value = -56874
$ws.Range("H114").Value = 37784.6
$ws.Range("J114").Value = 37784.6
$ws.Range("L114").Value = 37784.6
$ws.Range("N114").Value = -46462.6
$ws.Range("H116").Value = 49684
$ws.Range("J116").Value = 49684
$ws.Range("L116").Value = 49684
$ws.Range("N116").Value = -58862
$ws.Range("H130").Value = 50034.668
$ws.Range("J130").Value = 50034.668
$ws.Range("L130").Value = 50034.668
$ws.Range("N130").Value = -60074.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2813.125
$ws.Range("I7").Value = 3250
$ws.Range("J7").Value = 2667.5
$ws.Range("K7").Value = 3250
$ws.Range("L7").Value = 2667.5
$ws.Range("M7").Value = -3138
$ws.Range("N7").Value = -2891.5
$ws.Range("H22").Value = 604
$ws.Range("I22").Value = 604
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 604
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -309
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 604
$ws.Range("I27").Value = 604
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 604
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -497
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 2521.625
$ws.Range("I40").Value = 2422.4546
$ws.Range("K40").Value = 2422.4546
$ws.Range("M40").Value = -2286.4546
$ws.Range("H108").Value = 48618
$ws.Range("J108").Value = 48618
$ws.Range("L108").Value = 48618
$ws.Range("N108").Value = -56298
$ws.Range("H112").Value = 28598
$ws.Range("J112").Value = 28598
$ws.Range("L112").Value = 28598
$ws.Range("N112").Value = -31552
$ws.Range("H116").Value = 36600
$ws.Range("J116").Value = 36600
$ws.Range("L116").Value = 36600
$ws.Range("N116").Value = -45778
$ws.Range("H118").Value = 43405
$ws.Range("J118").Value = 43405
$ws.Range("L118").Value = 43405
$ws.Range("N118").Value = -46719
$ws.Range("H119").Value = 43960
$ws.Range("J119").Value = 43960
$ws.Range("L119").Value = 43960
$ws.Range("N119").Value = -53636
$ws.Range("H120").Value = 53840
$ws.Range("J120").Value = 53840
$ws.Range("L120").Value = 53840
$ws.Range("N120").Value = -63516
$ws.Range("H124").Value = 39924.75
$ws.Range("J124").Value = 39924.75
$ws.Range("L124").Value = 39924.75
$ws.Range("N124").Value = -49744.75
$ws.Range("H126").Value = 2813.125
$ws.Range("I126").Value = 3250
$ws.Range("J126").Value = 2667.5
$ws.Range("K126").Value = 9750
$ws.Range("L126").Value = 8002.5
$ws.Range("M126").Value = -7280
$ws.Range("N126").Value = -12942.5
$ws.Range("H127").Value = 49797.5
$ws.Range("J127").Value = 49797.5
$ws.Range("L127").Value = 49797.5
$ws.Range("N127").Value = -59717.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H108").Value = 48618
$ws.Range("J108").Value = 48618
$ws.Range("L108").Value = 48618
$ws.Range("N108").Value = -56298
$ws.Range("H110").Value = 49868.668
$ws.Range("J110").Value = 49868.668
$ws.Range("L110").Value = 49868.668
$ws.Range("N110").Value = -58048.668
$ws.Range("H116").Value = 48686
$ws.Range("J116").Value = 48686
$ws.Range("L116").Value = 48686
$ws.Range("N116").Value = -57864
$ws.Range("H117").Value = 49309
$ws.Range("J117").Value = 49309
$ws.Range("L117").Value = 49309
$ws.Range("N117").Value = -58487
$ws.Range("H121").Value = 43886.5
$ws.Range("J121").Value = 43886.5
$ws.Range("L121").Value = 43886.5
$ws.Range("N121").Value = -47380.5
$ws.Range("H126").Value = 4202609.5
$ws.Range("I126").Value = 9804423
$ws.Range("J126").Value = 1250
$ws.Range("K126").Value = 29413269
$ws.Range("L126").Value = 3750
$ws.Range("M126").Value = -29410799
$ws.Range("N126").Value = -8690
$ws.Range("H128").Value = 50711
$ws.Range("J128").Value = 50711
$ws.Range("L128").Value = 50711
$ws.Range("N128").Value = -60671
